# Update the "想去人数" (F column) counts that changed between the two
# data pulls, on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2-6 hold the F-column values.
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1040
$wsExpo.Range("F3").Value = 31
$wsExpo.Range("F4").Value = 2246
$wsExpo.Range("F5").Value = 16
$wsExpo.Range("F6").Value = 494

# Sheet "全部类型": same events, but shifted down two rows (rows 4-8).
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1040
$wsAll.Range("F5").Value = 31
$wsAll.Range("F6").Value = 2246
$wsAll.Range("F7").Value = 16
$wsAll.Range("F8").Value = 494
